$d = $word.ActiveDocument

# 1) Typo fix: "incluído" -> "incluido"
#    (El contenido que provee Samsung en el Sitio web o Móvil (incluído, ...)
$d.Content.Find.Execute("incluído, entre otros", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "incluido, entre otros", 2)

# 2) Typo fix: "No prodrá utilizar" -> "No podrá utilizar"
$d.Content.Find.Execute("No prodrá utilizar", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "No podrá utilizar", 2)

# 3) Wording fix: "entidad ni (d) para interferir" -> "entidad, o (d) para interferir"
$d.Content.Find.Execute("entidad ni (d) para interferir", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "entidad, o (d) para interferir", 2)

# 4) Word choice fix: "proporciona a tu filial." -> "proporciona a tu subsidiaria."
$d.Content.Find.Execute("proporciona a tu filial.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "proporciona a tu subsidiaria.", 2)

# 5) Move the "_GoBack" bookmark from the start of the document (paragraph 1)
#    to the end of the last paragraph, right after the "Fecha [22.01.2025]" run.
#    A direct zero-width Range placed exactly at the end of the document body
#    mis-resolves in this host, so we temporarily tag the spot with a unique
#    marker, locate it with Find (non-zero-width match), collapse onto its
#    leading edge, drop the bookmark there, then strip the marker back out.
$d.Content.Find.Execute("Fecha [22.01.2025]", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Fecha [22.01.2025]@@BMARK@@", 2)

$found = $d.Content
$found.Find.Execute("@@BMARK@@", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0)
$found.Collapse(1)
$d.Bookmarks.Add("_GoBack", $found)

$d.Content.Find.Execute("@@BMARK@@", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2)
